$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sku = "298044"
$desc = "Диск DVD+R 10шт Mirex 4.7Gb 16x Cake box printable inkjet (UL130029A1L)"

$qty = @(8, 1, 8, 1, 8, 1)
$price = @(161, 161, 161, 161, 161, 161)

# Columns A and B hold text values (the SKU code and the product
# description), so force a text number-format before writing them --
# otherwise "298044" would be auto-coerced into a numeric cell.
$ws.Range("A1:B6").NumberFormat = "@"

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $sku
    $ws.Cells.Item($r, 2).Value = $desc
    $ws.Cells.Item($r, 3).Value = $qty[$i]
    $ws.Cells.Item($r, 4).Value = $price[$i]
}
